$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Malik Beasley", "SG,SF", "Detroit Pistons"),
    @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
    @("Patrick Williams", "PF", "Chicago Bulls"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
    @("Kelly Olynyk", "C", "Toronto Raptors"),
    @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Kel'el Ware", "PF,C", "Miami Heat"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Josh Hart", "SG,SF,PF", "New York Knicks"),
    @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves"),
    @("Alperen Sengün", "C", "Houston Rockets")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
